# Apply the changes described by the commit:
#  - add 8 new "Battery discharger" rows (168-175) to the unitdata sheet
#  - reconfigure the AutoFilter (Generator_ID column) to show
#    "Hydrogen processor" / "Hydrogen storage dimensioner" instead of
#    "Battery charger" / "Battery discharger" (extending it over the
#    new rows too)
#  - move the sheet selection to F48
#  - (attempted) reposition the saved window

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("unitdata")

# --- 1. Append the new "Battery discharger" rows ------------------------
$newRows = @(
    @(168, "FI00", "Battery discharger", "National Trends",    2025, 40),
    @(169, "DE00", "Battery discharger", "National Trends",    2025, 600),
    @(170, "SE04", "Battery discharger", "National Trends",    2025, 70),
    @(171, "FR00", "Battery discharger", "National Trends",    2025, 900),
    @(172, "UK00", "Battery discharger", "National Trends",    2025, 1200),
    @(173, "PL00", "Battery discharger", "Distributed Energy", 2040, 5000),
    @(174, "DE00", "Battery discharger", "Distributed Energy", 2040, 25000),
    @(175, "NL00", "Battery discharger", "Distributed Energy", 2040, 5000)
)

foreach ($r in $newRows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
}

# --- 2. Re-apply the AutoFilter with the new criteria & range -----------
# Turn the existing filter off first so the filter range can grow to
# include the freshly-added rows.
$ws.AutoFilterMode = $false

[void]$ws.Range("A1:F175").AutoFilter(2, @("Hydrogen processor", "Hydrogen storage dimensioner"), 7)

# Keep the _xlnm._FilterDatabase defined name in sync with the new range.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=unitdata!`$A`$1:`$F`$175"
    }
}

# --- 3. Update the active selection on the sheet -------------------------
[void]$ws.Range("F48").Select()

# --- 4. Reposition the saved window (best effort) ------------------------
$win = $wb.Windows.Item(1)
$win.Left = -120
$win.Top = -120
